# Update "想去人数" (want-to-go count) / occasional "最低票价" (min price) figures
# across the four sheets, and append one new row to 本地生活 (Local Life) for a
# newly-scraped event, matching the upstream gh-pages data refresh @ 456a3b4.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: 展览 (Exhibitions)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1820
$ws.Range("F5").Value = 33
$ws.Range("F6").Value = 1068
$ws.Range("F8").Value = 162
$ws.Range("F9").Value = 555
$ws.Range("F10").Value = 45
$ws.Range("F11").Value = 437
$ws.Range("F12").Value = 186
$ws.Range("F13").Value = 1325
$ws.Range("F15").Value = 1374
$ws.Range("F17").Value = 53
$ws.Range("F18").Value = 269
$ws.Range("F19").Value = 1525
$ws.Range("F22").Value = 299
$ws.Range("F25").Value = 1114
$ws.Range("F26").Value = 301
$ws.Range("F27").Value = 779
$ws.Range("F29").Value = 971
$ws.Range("F30").Value = 199307
$ws.Range("F31").Value = 919
$ws.Range("F33").Value = 1307
$ws.Range("F34").Value = 871
$ws.Range("F36").Value = 10
$ws.Range("F37").Value = 810
$ws.Range("F38").Value = 1524
$ws.Range("F40").Value = 10
$ws.Range("F43").Value = 754
$ws.Range("F45").Value = 27

# ---------------------------------------------------------------------------
# Sheet: 演出 (Performances)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 104
$ws.Range("F11").Value = 1355
$ws.Range("F13").Value = 2440
$ws.Range("F14").Value = 1158
$ws.Range("F15").Value = 380
$ws.Range("F16").Value = 709
$ws.Range("F17").Value = 189
$ws.Range("F22").Value = 413
$ws.Range("F25").Value = 259
$ws.Range("F26").Value = 19698
$ws.Range("F31").Value = 225
$ws.Range("F33").Value = 39
$ws.Range("G34").Value = 266
$ws.Range("F35").Value = 6
$ws.Range("F38").Value = 156
$ws.Range("F42").Value = 15
$ws.Range("F43").Value = 15
$ws.Range("F45").Value = 108
$ws.Range("F46").Value = 49

# ---------------------------------------------------------------------------
# Sheet: 本地生活 (Local Life)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 858
$ws.Range("F6").Value = 2659
$ws.Range("F7").Value = 4426
$ws.Range("F10").Value = 493
$ws.Range("F11").Value = 546
$ws.Range("F12").Value = 377
$ws.Range("F13").Value = 52
$ws.Range("F14").Value = 493
$ws.Range("F15").Value = 160

# Append new row 16 - copy formatting of the A-column id cell (bold/bordered/
# centered style) from the row above so the new id cell matches the existing
# look, then fill in all the values for the new event.
$ws.Cells.Item(15, 1).Copy($ws.Cells.Item(16, 1))
$ws.Range("A16").Value = 15
# Force the date column to stay plain text (matches every other row, which
# stores "2024-07-05" as a string, not an Excel date serial).
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "2024-07-05"
$ws.Range("B16").ClearFormats()
$ws.Range("C16").Value = "上海·「剧场版BLUE LOCK -EPISODE凪- 」"
$ws.Range("D16").Value = "西藏北路198号大悦城北座8楼N809-1 animate cafe上海店"
$ws.Range("E16").Value = "2024.07.05 00:00-07.30 23:59"
$ws.Range("F16").Value = 64
$ws.Range("G16").Value = 30
$ws.Range("H16").Value = "https://show.bilibili.com/platform/detail.html?id=87171"
$ws.Range("I16").Value = "//i1.hdslb.com/bfs/openplatform/202406/9TBGVLTf1718097565516.png"

# ---------------------------------------------------------------------------
# Sheet: 全部类型 (All Types)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1820
$ws.Range("F3").Value = 858
$ws.Range("F5").Value = 4426
$ws.Range("F7").Value = 546
$ws.Range("F8").Value = 33
$ws.Range("F9").Value = 52
$ws.Range("F10").Value = 52
$ws.Range("F11").Value = 493
$ws.Range("F12").Value = 160
$ws.Range("F14").Value = 1068
$ws.Range("F16").Value = 162
$ws.Range("F18").Value = 1355
$ws.Range("F19").Value = 555
$ws.Range("F20").Value = 437
$ws.Range("F21").Value = 186
$ws.Range("F22").Value = 2440
$ws.Range("F23").Value = 1158
$ws.Range("F24").Value = 1325
$ws.Range("F26").Value = 1374
$ws.Range("F27").Value = 53
$ws.Range("F28").Value = 189
$ws.Range("F30").Value = 1525
$ws.Range("F32").Value = 299
$ws.Range("F33").Value = 413
$ws.Range("F34").Value = 1114
$ws.Range("F35").Value = 779
$ws.Range("F37").Value = 971
$ws.Range("F38").Value = 259
$ws.Range("F39").Value = 919
$ws.Range("F40").Value = 871
$ws.Range("F41").Value = 810
$ws.Range("G42").Value = 266
$ws.Range("F43").Value = 1524
$ws.Range("F45").Value = 156
$ws.Range("F48").Value = 15
$ws.Range("F49").Value = 754
$ws.Range("F51").Value = 27
